# spring 24 week 3 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.23
$ws.Range("E3").Value = 1.31
$ws.Range("G3").Value = 0.65
$ws.Range("B4").Value = 1.52
$ws.Range("E4").Value = 1.23
$ws.Range("G4").Value = 0.96
$ws.Range("C5").Value = 1.34
$ws.Range("D5").Value = 1.34
$ws.Range("F5").Value = 1.03
$ws.Range("E6").Value = 1.33
$ws.Range("F6").Value = 1.19
$ws.Range("C7").Value = 2.15
$ws.Range("D7").Value = 1.76
